$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing summary rows (28 and 29) so the table ends at row 27
$ws.Rows("28:29").Delete()

$rows = @(
    @(1661788527, 1, 0.1, 0.1, 1),
    @(1902183089, 2, 3, 3, 1),
    @(1906344657, 1, 3, 3, 1),
    @(1500015628, 1, 0.1, 0.55, 1),
    @(1522032606, 3, 2.1, 2.1, 1),
    @(1545108803, 2, 0.1, 0.1, 1),
    @(830471345, 0.5, 0.1, 0.1, 1),
    @(100696801, 0.1, 0.1, 0.1, 1),
    @(1100805552, 0.1, 3, 2.47, 1),
    @(1214250561, 2, 0.98, 0.63, 1),
    @(1371274123, 0.1, 2.48, 3, 1),
    @(1610612861, 0.1, 0.1, 0.1, 1),
    @(1809382033, 0.1, 0.1, 0.1, 1),
    @(218891067, 2, 1.85, 0.1, 1),
    @(26278542, 0.5, 0.1, 0.1, 1),
    @(870843266, 2, 3, 3, 1),
    @(932970289, 0.5, 0.1, 0.1, 1),
    @(1173749970, 3, 0.1, 0.1, 1),
    @(1251868241, 0.3, 2.07, 1.72, 1),
    @(1333067343, 1, 0.35, 1.15, 1),
    @(1902706256, 0.1, 0.47, 0.1, 1),
    @(829424033, 0.5, 1.32, 3, 1),
    @(973684048, 0.1, 1.92, 0.1, 1),
    @("RF", 14, 47, 37, 12),
    @("IRF", 806, 1763, 2272, 679),
    @("MOSD", 0.301, 0.351, 0.36, 0.299)
)

# Write the refreshed portfolio data into A2:E27
for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}

# Match the saved selection/active range to the new used range
[void]$ws.Range("A1:E27").Select()
